$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (row 1 header text)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 12:55"

# Full data refresh: province/city rows 4-64 re-sorted by "Casos totales"
# (descending) with updated counts for several provinces.
$rows = @(
    @{ Row=4; A="Madrid"; B=24090; C=8301; D=12397; E=3392 },
    @{ Row=5; A="Cataluña"; B=15026; C=3455; D=10345; E=1226 },
    @{ Row=6; A="Castilla-La Mancha"; B=2780; C=71; D=2446; E=263 },
    @{ Row=7; A="Bizkaia/Vizcaya"; B=2776; C=1503; D=2010; E=116 },
    @{ Row=8; A="Valencia/Valencia"; B=2685; C=130; D=2412; E=143 },
    @{ Row=9; A="Navarra"; B=2146; C=161; D=1883; E=102 },
    @{ Row=10; A="Araba/Alava"; B=1947; C=1503; D=1332; E=115 },
    @{ Row=11; A="Alacant/Alicante"; B=1839; C=46; D=1658; E=135 },
    @{ Row=12; A="Ciudad Real"; B=1755; C=236; D=1525; E=145 },
    @{ Row=13; A="Zaragoza"; B=1641; C=141; D=1419; E=81 },
    @{ Row=14; A="La Rioja"; B=1629; C=397; D=1164; E=68 },
    @{ Row=15; A="Albacete"; B=1386; C=236; D=1204; E=122 },
    @{ Row=16; A="A Coruña"; B=1351; C=153; D=1261; E=40 },
    @{ Row=17; A="Malaga"; B=1321; C=83; D=1169; E=69 },
    @{ Row=18; A="Toledo"; B=1317; C=236; D=1126; E=131 },
    @{ Row=19; A="Asturias"; B=1158; C=78; D=1032; E=48 },
    @{ Row=20; A="Cantabria"; B=1100; C=24; D=1049; E=27 },
    @{ Row=21; A="Pontevedra"; B=1060; C=153; D=1005; E=9 },
    @{ Row=22; A="Sevilla"; B=1052; C=18; D=1000; E=34 },
    @{ Row=23; A="Caceres"; B=1045; C=11; D=945; E=89 },
    @{ Row=24; A="Salamanca"; B=1030; C=157; D=774; E=99 },
    @{ Row=25; A="Gipuzkoa/Guipuzcoa"; B=1017; C=1503; D=630; E=34 },
    @{ Row=26; A="Granada"; B=963; C=15; D=882; E=66 },
    @{ Row=27; A="Murcia"; B=939; C=17; D=897; E=25 },
    @{ Row=28; A="Aragon"; B=907; C=29; D=838; E=40 },
    @{ Row=29; A="Valladolid"; B=886; C=127; D=702; E=57 },
    @{ Row=30; A="Leon"; B=821; C=118; D=626; E=77 },
    @{ Row=31; A="Burgos"; B=719; C=137; D=533; E=49 },
    @{ Row=32; A="La Palma"; B=712; C=30; D=1056; E=2 },
    @{ Row=33; A="Jaen"; B=599; C=17; D=559; E=23 },
    @{ Row=34; A="Castello/Castellon"; B=586; C=9; D=545; E=32 },
    @{ Row=35; A="Cordoba"; B=572; C=4; D=555; E=13 },
    @{ Row=36; A="Segovia"; B=567; C=131; D=378; E=58 },
    @{ Row=37; A="Guadalajara"; B=535; C=236; D=436; E=86 },
    @{ Row=38; A="Soria"; B=523; C=61; D=432; E=30 },
    @{ Row=39; A="Badajoz"; B=515; C=49; D=449; E=17 },
    @{ Row=40; A="Cadiz"; B=507; C=10; D=484; E=13 },
    @{ Row=41; A="Ourense"; B=458; C=153; D=415; E=8 },
    @{ Row=42; A="Avila"; B=414; C=82; D=291; E=41 },
    @{ Row=43; A="Fuerteventura"; B=288; C=30; D=1056; E=0 },
    @{ Row=44; A="Lugo"; B=270; C=153; D=244; E=4 },
    @{ Row=45; A="Palencia"; B=262; C=28; D=221; E=13 },
    @{ Row=46; A="Cuenca"; B=253; C=236; D=180; E=55 },
    @{ Row=47; A="Almeria"; B=223; C=6; D=203; E=14 },
    @{ Row=48; A="Teruel"; B=222; C=14; D=196; E=12 },
    @{ Row=49; A="Huesca"; B=215; C=19; D=185; E=11 },
    @{ Row=50; A="Mallorca"; B=210; C=18; D=194; E=12 },
    @{ Row=51; A="Zamora"; B=192; C=30; D=144; E=18 },
    @{ Row=52; A="Huelva"; B=168; C=2; D=162; E=4 },
    @{ Row=53; A="Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"; B=58; C=0; D=58; E=3 },
    @{ Row=54; A="Melilla"; B=51; C=0; D=50; E=1 },
    @{ Row=55; A="Lanzarote"; B=49; C=30; D=1056; E=3 },
    @{ Row=56; A="La Gomera"; B=35; C=30; D=1056; E=0 },
    @{ Row=57; A="El Hierro"; B=30; C=30; D=1056; E=0 },
    @{ Row=58; A="Ceuta"; B=29; C=0; D=28; E=1 },
    @{ Row=59; A="Ibiza"; B=21; C=18; D=20; E=1 },
    @{ Row=60; A="Menorca"; B=15; C=18; D=13; E=0 },
    @{ Row=61; A="Gran Canaria"; B=8; C=30; D=1056; E=11 },
    @{ Row=62; A="Arroyo de la Luz"; B=7; C=0; D=7; E=0 },
    @{ Row=63; A="Tenerife"; B=3; C=30; D=1056; E=36 },
    @{ Row=64; A="Formentera"; B=0; C=10; D=0; E=8 }

)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
